$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly price update: a new record is inserted at row 49, pushing the
# existing rows 49-56 down to 50-57 (each retains its original data).
$ws.Rows.Item(49).Insert()

# Populate the new row 49 with this week's data.
$ws.Range("A49").Value = 7
$ws.Range("B49").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C49").Value = "Ñuble"
$ws.Range("D49").Value = 44875
$ws.Range("E49").Value = 16
$ws.Range("F49").Value = 100112001
$ws.Range("G49").Value = "Berenjena"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 60
$ws.Range("K49").Value = 10000
$ws.Range("L49").Value = 11000
$ws.Range("M49").Value = 10500
$ws.Range("N49").Value = "$/caja 60 unidades"
$ws.Range("O49").Value = "Región de Arica y Parinacota"
$ws.Range("P49").Value = 175
$ws.Range("Q49").Value = 60
$ws.Range("R49").Value = "Hortaliza"
